$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had two separate MDR-TB introduction parameters:
#   row 5: start_mdr_introduce_time | 1880 | "Calendar year that MDR-TB first begins to emerge"
#   row 6: end_mdr_introduce_time   | 1885 | "Calendar year that MDR-TB amplification reaches full parameter value"
#
# Per the commit message, drop the "end" row entirely and rename the
# remaining "start" parameter to simply "mdr_introduce_time" (keeping its
# original value/description).

# Delete the entire "end_mdr_introduce_time" row (row 6); rows below shift up.
$ws.Rows("6").Delete() | Out-Null

# Rename the now-single MDR introduction parameter (originally row 5,
# "start_mdr_introduce_time") to "mdr_introduce_time".
$ws.Range("A5").Value = "mdr_introduce_time"

# Match the final selection left behind in the saved workbook.
$ws.Range("C6:D8").Select() | Out-Null
